# Updated cryptos list on Sat Oct  5 21:28:22 UTC 2024 with GitHub Actions
#
# Refreshes the hourly snapshot of coin Price (col D) / Volume(1h) (col E)
# pulled from coinranking.com, and reflects that "RenderToken" overtook
# "Monero" in the rankings, so rows 37/38 swap their Coin/Link/Price/Volume
# data.
#
# All of these cells are plain text in the workbook (prices are stored as
# formatted strings like "61.891.85" or "0.0000172", not numbers, and the
# volume column holds padded strings like "  -0.88%  "). Excel's COM layer
# auto-coerces an ambiguous numeric-looking string typed into .Value into a
# real number (e.g. "1.00" -> 1, dropping the trailing zero), so each cell
# is forced to Text format before the write, then restored to the workbook's
# normal (unstyled) look afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$updates = @(
    @{ Cell = "D2";  Value = "61.891.85" },
    @{ Cell = "E2";  Value = "  -0.88%  " },
    @{ Cell = "D3";  Value = "2.408.54" },
    @{ Cell = "E3";  Value = "  -0.77%  " },
    @{ Cell = "D4";  Value = "1.00" },
    @{ Cell = "E4";  Value = "  -0.06%  " },
    @{ Cell = "D5";  Value = "561.34" },
    @{ Cell = "E5";  Value = "  +0.82%  " },
    @{ Cell = "D6";  Value = "142.14" },
    @{ Cell = "E6";  Value = "  -1.19%  " },
    @{ Cell = "E7";  Value = "  -0.02%  " },
    @{ Cell = "D8";  Value = "0.530" },
    @{ Cell = "E8";  Value = "  -0.97%  " },
    @{ Cell = "D9";  Value = "0.109" },
    @{ Cell = "E9";  Value = "  -1.47%  " },
    @{ Cell = "E10"; Value = "  -1.86%  " },
    @{ Cell = "D11"; Value = "5.24" },
    @{ Cell = "E11"; Value = "  -3.18%  " },
    @{ Cell = "D12"; Value = "0.348" },
    @{ Cell = "E12"; Value = "  -0.86%  " },
    @{ Cell = "D13"; Value = "25.46" },
    @{ Cell = "E13"; Value = "  -3.45%  " },
    @{ Cell = "D14"; Value = "0.0000172" },
    @{ Cell = "E14"; Value = "  -2.18%  " },
    @{ Cell = "D15"; Value = "2.839.20" },
    @{ Cell = "E15"; Value = "  -0.91%  " },
    @{ Cell = "D16"; Value = "61.911.30" },
    @{ Cell = "E16"; Value = "  -0.71%  " },
    @{ Cell = "D17"; Value = "2.406.93" },
    @{ Cell = "E17"; Value = "  -0.85%  " },
    @{ Cell = "D18"; Value = "11.19" },
    @{ Cell = "E18"; Value = "  +0.76%  " },
    @{ Cell = "D19"; Value = "320.71" },
    @{ Cell = "E19"; Value = "  -1.31%  " },
    @{ Cell = "D20"; Value = "6.80" },
    @{ Cell = "E20"; Value = "  +0.73%  " },
    @{ Cell = "D21"; Value = "4.12" },
    @{ Cell = "E21"; Value = "  -1.71%  " },
    @{ Cell = "E22"; Value = "  -0.40%  " },
    @{ Cell = "D23"; Value = "65.31" },
    @{ Cell = "E23"; Value = "  +0.47%  " },
    @{ Cell = "E24"; Value = "  -3.59%  " },
    @{ Cell = "D25"; Value = "8.68" },
    @{ Cell = "E25"; Value = "  -4.57%  " },
    @{ Cell = "D26"; Value = "561.87" },
    @{ Cell = "E26"; Value = "  -2.45%  " },
    @{ Cell = "D27"; Value = "0.995" },
    @{ Cell = "E27"; Value = "  -0.51%  " },
    @{ Cell = "D28"; Value = "2.515.35" },
    @{ Cell = "E28"; Value = "  -1.22%  " },
    @{ Cell = "D29"; Value = "0.0₃0931" },
    @{ Cell = "E29"; Value = "  -1.81%  " },
    @{ Cell = "D30"; Value = "8.14" },
    @{ Cell = "E30"; Value = "  -3.28%  " },
    @{ Cell = "D31"; Value = "1.38" },
    @{ Cell = "E31"; Value = "  -5.07%  " },
    @{ Cell = "D32"; Value = "0.147" },
    @{ Cell = "E32"; Value = "  -1.28%  " },
    @{ Cell = "E33"; Value = "  +0.13%  " },
    @{ Cell = "D34"; Value = "1.50" },
    @{ Cell = "E34"; Value = "  -4.72%  " },
    @{ Cell = "E35"; Value = "  +0.02%  " },
    @{ Cell = "D36"; Value = "4.75" },
    @{ Cell = "E36"; Value = "  -2.20%  " },

    # Row 37/38 swap: RenderToken overtakes Monero.
    @{ Cell = "B37"; Value = "RenderToken" },
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render" },
    @{ Cell = "D37"; Value = "5.43" },
    @{ Cell = "E37"; Value = "  -5.39%  " },
    @{ Cell = "B38"; Value = "Monero" },
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" },
    @{ Cell = "D38"; Value = "152.32" },
    @{ Cell = "E38"; Value = "  +1.18%  " },

    @{ Cell = "D39"; Value = "0.379" },
    @{ Cell = "E39"; Value = "  -1.67%  " },
    @{ Cell = "D40"; Value = "18.47" },
    @{ Cell = "E40"; Value = "  -1.76%  " },
    @{ Cell = "D41"; Value = "1.79" },
    @{ Cell = "E41"; Value = "  -5.45%  " },
    @{ Cell = "D43"; Value = "147.79" },
    @{ Cell = "E43"; Value = "  -2.24%  " },
    @{ Cell = "D44"; Value = "2.23" },
    @{ Cell = "E44"; Value = "  -5.08%  " },
    @{ Cell = "D45"; Value = "3.58" },
    @{ Cell = "E45"; Value = "  -1.70%  " },
    @{ Cell = "D46"; Value = "0.0526" },
    @{ Cell = "E46"; Value = "  -3.39%  " },
    @{ Cell = "D47"; Value = "19.83" },
    @{ Cell = "E47"; Value = "  -3.25%  " },
    @{ Cell = "D48"; Value = "0.590" },
    @{ Cell = "E48"; Value = "  -0.22%  " },
    @{ Cell = "D49"; Value = "0.0915" },
    @{ Cell = "E49"; Value = "  -0.48%  " },
    @{ Cell = "D50"; Value = "0.0225" },
    @{ Cell = "E50"; Value = "  -1.74%  " },
    @{ Cell = "D51"; Value = "11.53" },
    @{ Cell = "E51"; Value = "  +0.25%  " }
)

foreach ($u in $updates) {
    Set-TextValue $u.Cell $u.Value
}
